# "mise en place du systeme de frames"
#
# Updates the user-story rows describing alien movement (row 10) and
# player-laser firing (row 18) to reflect the new frame-based system, and
# moves the view/selection to the area that was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 10 : "Deplacement des aliens" ---------------------------------
# The two "moving right / moving left" test cases lose their trailing
# "Pour faires bouger les aliens" justification line, and the
# "reach the edge of the screen" test case is reworded. The fourth test
# case (the horizontal walk-cycle animation) is removed entirely.
$ws.Range("C10").Value = "Quand la partie est en cours`nLes aliens se déplacent vers la droite a un rythme régulier"
$ws.Range("D10").Value = "Quand la partie est en cours`nLes aliens se déplacent vers la gauche a un rythme régulier`n"
$ws.Range("E10").Value = "Quand un alien touche le bord de l'écran`nLes aliens descende d'une ligne et reparte dans l'autre sens`n"
$ws.Range("F10").Value = ""

# A10 becomes a "header" style row (like A11, A14, A17, A20) - copy the
# formatting only from A11 so the shared cell style is reused.
$ws.Range("A11").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 18 : "Tir du vaisseau du joueur" -------------------------------
# The laser-firing / rate-limiting test cases are reworded for the new
# missile-frame system, and two brand new test cases are appended
# describing the missile's upward travel and its disappearance.
$ws.Range("C18").Value = "En jeu`nQuand la flèche haut est préssée en jeu`nLe vaisseau tir un laser vers le haut`n"
$ws.Range("D18").Value = "Le laser a été tirer il y a moins d'une seconde`nQuand la flèche haut est préssée en jeu`nil ne se passe rien`n"
$ws.Range("E18").Value = "Le missile est tiré`nIl avance d'une case toute les quart de seconde vers le haut"
$ws.Range("F18").Value = "Le missile atteint le haut de l'écran`nil disparait"

# A18 also becomes a "header" style row.
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 18 grew an extra line of content, so it needs to be a little taller.
$ws.Rows.Item(18).RowHeight = 94.5

# --- View / selection ----------------------------------------------------
# The author ended up scrolled down to the rows they just edited, with
# B18 selected.
$ws.Activate()
$ws.Range("B18").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
